$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 908.1667
$ws.Range("I15").Value = 908.1667
$ws.Range("K15").Value = 2724.5001
$ws.Range("M15").Value = -2555.5001

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 569.2105
$ws.Range("I28").Value = 421.33334
$ws.Range("J28").Value = 1123.75
$ws.Range("K28").Value = 421.33334
$ws.Range("L28").Value = 1123.75
$ws.Range("M28").Value = 63.66665999999998
$ws.Range("N28").Value = -2093.75

# Row 40: Stuck in the Moment
$ws.Range("H40").Value = 2995
$ws.Range("I40").Value = 2984
$ws.Range("K40").Value = 2984
$ws.Range("M40").Value = -2809

# Row 51: A Bile Business
$ws.Range("H51").Value = 24497
$ws.Range("I51").Value = 52287
$ws.Range("K51").Value = 52287
$ws.Range("M51").Value = -51803

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 23836.055
$ws.Range("I62").Value = 17492.2
$ws.Range("K62").Value = 17492.2
$ws.Range("M62").Value = -16868.2

# Row 64: Forged from the Void
$ws.Range("H64").Value = 63360.223
$ws.Range("I64").Value = 7075.25
$ws.Range("J64").Value = 79441.64
$ws.Range("K64").Value = 7075.25
$ws.Range("L64").Value = 79441.64
$ws.Range("M64").Value = -6827.25
$ws.Range("N64").Value = -79937.64

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 23836.055
$ws.Range("I65").Value = 17492.2
$ws.Range("K65").Value = 87461
$ws.Range("M65").Value = -84341

# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 63360.223
$ws.Range("I67").Value = 7075.25
$ws.Range("J67").Value = 79441.64
$ws.Range("K67").Value = 7075.25
$ws.Range("L67").Value = 79441.64
$ws.Range("M67").Value = -6217.25
$ws.Range("N67").Value = -81157.64

# Row 113: Amaro Kart
$ws.Range("H113").Value = 6518.25
$ws.Range("I113").Value = 3999.75
$ws.Range("J113").Value = 7777.5
$ws.Range("K113").Value = 3999.75
$ws.Range("L113").Value = 7777.5
$ws.Range("M113").Value = -745.75
$ws.Range("N113").Value = -14285.5

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 419729.44
$ws.Range("I132").Value = 450690.47
$ws.Range("K132").Value = 1352071.41
$ws.Range("M132").Value = -1349541.41

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1208.2222
$ws.Range("I141").Value = 1208.2222
$ws.Range("K141").Value = 3624.6666
$ws.Range("M141").Value = 1555.3334

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 6256271.5
$ws.Range("I32").Value = 7250593
$ws.Range("K32").Value = 7250593
$ws.Range("M32").Value = -7250306

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 3332.7727
$ws.Range("I45").Value = 3122.5264
$ws.Range("K45").Value = 3122.5264
$ws.Range("M45").Value = -2745.5264

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 6940.2256
$ws.Range("I61").Value = 3854.4348
$ws.Range("J61").Value = 15811.875
$ws.Range("K61").Value = 3854.4348
$ws.Range("L61").Value = 15811.875
$ws.Range("M61").Value = -3642.4348
$ws.Range("N61").Value = -16235.875

# Row 130: A Gift of Gloves
$ws.Range("H130").Value = 19997.5
$ws.Range("J130").Value = 19997.5
$ws.Range("L130").Value = 19997.5
$ws.Range("N130").Value = -30037.5

# Row 131: Additions to the Armoire
$ws.Range("H131").Value = 99996.664
$ws.Range("J131").Value = 99996.664
$ws.Range("L131").Value = 99996.664
$ws.Range("N131").Value = -110076.664

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 1016997.06
$ws.Range("I132").Value = 1296057.5
$ws.Range("K132").Value = 3888172.5
$ws.Range("M132").Value = -3885642.5

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 6940.2256
$ws.Range("I136").Value = 3854.4348
$ws.Range("J136").Value = 15811.875
$ws.Range("K136").Value = 11563.3044
$ws.Range("L136").Value = 47435.625
$ws.Range("M136").Value = -9013.304400000001
$ws.Range("N136").Value = -52535.625

$ws = $wb.Worksheets.Item("BSM")
# Row 96: Hammer Time
$ws.Range("H96").Value = 18433
$ws.Range("I96").Value = 18433
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 18433
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("M96").Value = -15687

# Row 133: Paring Is Caring
$ws.Range("H133").Value = 82541.336
$ws.Range("J133").Value = 90000
$ws.Range("L133").Value = 90000
$ws.Range("N133").Value = -100120

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 1385.5714
$ws.Range("I16").Value = 1385.5714
$ws.Range("K16").Value = 1385.5714
$ws.Range("M16").Value = -1098.5714

# Row 31: Wall Not Found
$ws.Range("H31").Value = 6877.387
$ws.Range("I31").Value = 1468.381
$ws.Range("J31").Value = 9647.853999999999
$ws.Range("K31").Value = 1468.381
$ws.Range("L31").Value = 9647.853999999999
$ws.Range("M31").Value = -1173.381
$ws.Range("N31").Value = -10237.854

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 6877.387
$ws.Range("I34").Value = 1468.381
$ws.Range("J34").Value = 9647.853999999999
$ws.Range("K34").Value = 1468.381
$ws.Range("L34").Value = 9647.853999999999
$ws.Range("M34").Value = -1266.381
$ws.Range("N34").Value = -10051.854

# Row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 1031336.75
$ws.Range("I58").Value = 1123731
$ws.Range("K58").Value = 1123731
$ws.Range("M58").Value = -1123528

# Row 74: License to Heal
$ws.Range("H74").Value = 44610.777
$ws.Range("J74").Value = 44647
$ws.Range("L74").Value = 44647
$ws.Range("N74").Value = -46395

# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 44610.777
$ws.Range("J77").Value = 44647
$ws.Range("L77").Value = 133941
$ws.Range("N77").Value = -142677

# Row 113: Patient Patients
$ws.Range("H113").Value = 1385.5714
$ws.Range("I113").Value = 1385.5714
$ws.Range("K113").Value = 1385.5714
$ws.Range("M113").Value = 784.4286

# Row 136: Turali Quality
$ws.Range("H136").Value = 1031336.75
$ws.Range("I136").Value = 1123731
$ws.Range("K136").Value = 3371193
$ws.Range("M136").Value = -3368643

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers
$ws.Range("H107").Value = 522.13336
$ws.Range("I107").Value = 304.22223
$ws.Range("K107").Value = 304.22223
$ws.Range("M107").Value = 1615.77777

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 3292.9285
$ws.Range("I93").Value = 2566.5557
$ws.Range("J93").Value = 4600.4
$ws.Range("K93").Value = 2566.5557
$ws.Range("L93").Value = 4600.4
$ws.Range("M93").Value = -1318.5557
$ws.Range("N93").Value = -7096.4

# Row 132: Tenets of Tanning
$ws.Range("H132").Value = 890909.3
$ws.Range("I132").Value = 1051165.6
$ws.Range("K132").Value = 3153496.8
$ws.Range("M132").Value = -3150966.8

# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 10042.538
$ws.Range("I136").Value = 9232.091
$ws.Range("K136").Value = 27696.273
$ws.Range("M136").Value = -25146.273

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Range("H2").Value = 83761650
$ws.Range("I2").Value = 12475.125
$ws.Range("K2").Value = 12475.125
$ws.Range("M2").Value = -12363.125

# Row 113: A Tender Table
$ws.Range("H113").Value = 2503.8462
$ws.Range("J113").Value = 2422.889
$ws.Range("L113").Value = 7268.667
$ws.Range("N113").Value = -11608.667

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2921.889
$ws.Range("I126").Value = 2916.8518
$ws.Range("J126").Value = 2937
$ws.Range("K126").Value = 8750.555399999999
$ws.Range("L126").Value = 8811
$ws.Range("M126").Value = -6280.555399999999
$ws.Range("N126").Value = -13751

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 21891592
$ws.Range("I132").Value = 2305433.5
$ws.Range("K132").Value = 6916300.5
$ws.Range("M132").Value = -6913770.5

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 10884764
$ws.Range("I136").Value = 11543765
$ws.Range("K136").Value = 34631295
$ws.Range("M136").Value = -34628745
